# Doing Updates for Financials
# Insert two new quarterly columns (D:E) in front of the existing quarterly
# data on the "LULU" sheet, shifting the historical quarters right, and
# populate the two new columns with the new quarter's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LULU")

# Insert 2 new columns at D:E; existing D:K shift right to F:M.
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (the old column D, now shifted) onto
# the two freshly inserted columns so they pick up the same cell styles
# (date format for header rows, number format for data rows) instead of
# falling back to the worksheet default.
for ($r = 5; $r -le 102; $r++) {
    $ws.Range("F" + $r).Copy()
    $ws.Range("D" + $r + ":E" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Row 7 / 38 / 80 - "Period Ending" header dates
$ws.Range("D7").Value = 43499
$ws.Range("E7").Value = 43401
$ws.Range("D38").Value = 43499
$ws.Range("E38").Value = 43401
$ws.Range("D80").Value = 43499
$ws.Range("E80").Value = 43401

# --- Income Statement (rows 8-35) ---
$ws.Range("D8").Value = 1167500
$ws.Range("E8").Value = 747700
$ws.Range("D9").Value = 498900
$ws.Range("E9").Value = 340900
$ws.Range("D10").Value = 668600
$ws.Range("E10").Value = 406800
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 836000
$ws.Range("E17").Value = 611800
$ws.Range("D18").Value = 331500
$ws.Range("E18").Value = 135900
$ws.Range("D20").Value = 2800
$ws.Range("E20").Value = 2000
$ws.Range("D21").Value = 369700
$ws.Range("E21").Value = 169600
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 334300
$ws.Range("E23").Value = 137900
$ws.Range("D24").Value = 113500
$ws.Range("E24").Value = 38400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 220800
$ws.Range("E26").Value = 99600
$ws.Range("D27").Value = 220800
$ws.Range("E27").Value = 99600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -2300
$ws.Range("E29").Value = -5200
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -2800
$ws.Range("E32").Value = -2000
$ws.Range("D33").Value = 218500
$ws.Range("E33").Value = 94400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 218500
$ws.Range("E35").Value = 94400

# --- Balance Sheet (rows 41-77) ---
$ws.Range("D41").Value = 881300
$ws.Range("E41").Value = 703600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 35800
$ws.Range("E43").Value = 29400
$ws.Range("D44").Value = 404800
$ws.Range("E44").Value = 496000
$ws.Range("D45").Value = 107300
$ws.Range("E45").Value = 134400
$ws.Range("D46").Value = 1429300
$ws.Range("E46").Value = 1363500
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 567200
$ws.Range("E48").Value = 531300
$ws.Range("D49").Value = 24200
$ws.Range("E49").Value = 24200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 64000
$ws.Range("E52").Value = 62100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2084700
$ws.Range("E54").Value = 1981000
$ws.Range("D57").Value = 95500
$ws.Range("E57").Value = 153100
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 404900
$ws.Range("E59").Value = 290500
$ws.Range("D60").Value = 500500
$ws.Range("E60").Value = 443700
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 138300
$ws.Range("E62").Value = 130600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 638700
$ws.Range("E66").Value = 574300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1346900
$ws.Range("E72").Value = 1310500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1446000
$ws.Range("E76").Value = 1406800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0

# --- Cash Flow Statement (rows 81-102) ---
$ws.Range("D81").Value = 218500
$ws.Range("E81").Value = 94400
$ws.Range("D83").Value = 35400
$ws.Range("E83").Value = 31700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 425900
$ws.Range("E89").Value = 106900
$ws.Range("D91").Value = -69100
$ws.Range("E91").Value = -72700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -76900
$ws.Range("E94").Value = -76600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -183900
$ws.Range("E100").Value = -106100
$ws.Range("D101").Value = 12500
$ws.Range("E101").Value = 1700
$ws.Range("D102").Value = 177700
$ws.Range("E102").Value = -74200
